$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-31 06:54:12"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-31 06:54:03"
$zhcn.Range("K2").Value = "2016-08-31 06:54:31"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-31 06:54:12"
$dede.Range("K2").Value = "2016-08-31 06:54:38"
